# edit.ps1 - apply the documented change set to the Word document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Using FIFO sequence ... into the array." -> "... into the array list."
#    (net textual effect of splitting that run into three runs)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "into the array. The order",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "into the array list. The order", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Eg. Array[0] = ..." -> "Eg. List[0] = ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Array[",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "List[", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Relocate the hidden "_GoBack" bookmark that currently sits between the
#    "Datatypes" run and the " in excel" run of the last paragraph. It needs
#    to end up (collapsed) at the end of a brand-new paragraph added further
#    down ("Where the key of the dictionary"). Deleting it here just merges
#    the two text runs around it back together - which is exactly the text
#    already present, so nothing else needs to change in this paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 4) Append the new paragraphs after the last paragraph ("Copy over Data
#    Types to Datatypes in excel"), all as plain (non-list) paragraphs.
# ---------------------------------------------------------------------------
function Add-PlainParagraph($afterRange, [string]$text) {
    $afterRange.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Last
    $newPara.Range.ListFormat.RemoveNumbers() | Out-Null
    $newPara.Style = "Normal"
    if ($text) {
        $newPara.Range.InsertAfter($text) | Out-Null
    }
    return $d.Paragraphs.Last.Range
}

$cursor = $d.Paragraphs.Last.Range
$cursor = Add-PlainParagraph $cursor ""
$cursor = Add-PlainParagraph $cursor "Notes: We are using Lists instead of arrays for the dynamic resizing feature of lists"
$cursor = Add-PlainParagraph $cursor ""
$cursor = Add-PlainParagraph $cursor "We have 2 Data structures. A Dictionary and a List."
$cursor = Add-PlainParagraph $cursor ""
$cursor = Add-PlainParagraph $cursor "Where the key of the dictionary"

# Remember this paragraph - the relocated bookmark goes at its very end.
$whereKeyParaIndex = $d.Paragraphs.Count

$cursor = Add-PlainParagraph $cursor ""
$cursor = Add-PlainParagraph $cursor ""
$cursor = Add-PlainParagraph $cursor "Changes to write to in Excel file."
$cursor = Add-PlainParagraph $cursor "A3 cell = name of the tab"
$cursor = Add-PlainParagraph $cursor ""

# ---------------------------------------------------------------------------
# 5) Put the "_GoBack" bookmark (collapsed) right after the text of the
#    "Where the key of the dictionary" paragraph, before its paragraph mark.
#    Placing a collapsed range/bookmark exactly on that boundary position is
#    mishandled by the host, so nudge it: temporarily extend the paragraph by
#    one throw-away character, anchor the bookmark just before that
#    character (now safely off the boundary), then remove the throw-away
#    character again - the bookmark stays put.
# ---------------------------------------------------------------------------
$whereKeyPara = $d.Paragraphs.Item($whereKeyParaIndex)
$targetPos = $whereKeyPara.Range.End - 1

$whereKeyPara.Range.InsertAfter("X") | Out-Null

$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$whereKeyPara2 = $d.Paragraphs.Item($whereKeyParaIndex)
$throwAwayPos = $whereKeyPara2.Range.End - 2
$throwAwayRange = $d.Range($throwAwayPos, $throwAwayPos + 1)
$throwAwayRange.Delete() | Out-Null

Write-Host "Final paragraph dump:"
foreach ($p in $d.Paragraphs) {
    Write-Host "---PARA---"
    Write-Host $p.Range.Text
}
